$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.681.89"
$ws.Range("D3").Value = "1.643.31"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.92"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.07"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "1.871.70"

# Rows 13 and 14 swap content (Polkadot <-> WrappedEther) with value updates
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.634.54"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").Value = "26.696.49"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.30"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  +14.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.43"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +4.60%  "
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0514"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").Value = "1.280.23"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.532"
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("E39").Value = "  +3.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").Value = "1.781.83"
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.84"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.20"
$ws.Range("E46").Value = "  +7.83%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.77"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("E51").Value = "  -0.47%  "
